$wb = $excel.ActiveWorkbook

# ---- Sheet 1: LP1912 ----
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("A2").Value = "Última actualización: 03:02:37"
$ws1.Range("A3").Value = "Total filas: 4"

# Row 6
$ws1.Range("A6").Value = "03:02:37"
$ws1.Range("B6").Value = "03:48"
$ws1.Range("C6").Value = "14_ABASTO"
$ws1.Range("D6").Value = 46
$ws1.Range("E6").Value = "LP1912"

# Row 7
$ws1.Range("A7").Value = "03:02:37"
$ws1.Range("B7").Value = "04:01"
$ws1.Range("C7").Value = "81_EL PELIGRO"
$ws1.Range("D7").Value = 59
$ws1.Range("E7").Value = "LP1912"

# Row 8
$ws1.Range("A8").Value = "03:02:37"
$ws1.Range("B8").Value = "04:45"
$ws1.Range("C8").Value = "215A_EL PATO"
$ws1.Range("D8").Value = 103
$ws1.Range("E8").Value = "LP1912"

# Row 9 (new)
$ws1.Range("A9").Value = "03:02:37"
$ws1.Range("B9").Value = "04:53"
$ws1.Range("C9").Value = "11_ETCHEVERRY"
$ws1.Range("D9").Value = 111
$ws1.Range("E9").Value = "LP1912"

# ---- Sheet 2: LP1912-215 ----
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("A2").Value = "Última actualización: 03:02:37"

# Row 6
$ws2.Range("A6").Value = "03:02:37"
$ws2.Range("B6").Value = "04:45"
$ws2.Range("C6").Value = "215A_EL PATO"
$ws2.Range("D6").Value = 103
$ws2.Range("E6").Value = "LP1912"

# ---- Sheet 3: 6203-6173 ----
$ws3 = $wb.Worksheets.Item(3)

$ws3.Range("A2").Value = "Última actualización: 03:02:37"
